$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 9267.571
$ws.Range("I5").Value = 8425.6
$ws.Range("K5").Value = 8425.6
$ws.Range("M5").Value = -8310.6
$ws.Range("H51").Value = 8278.296
$ws.Range("J51").Value = 9098.286
$ws.Range("L51").Value = 9098.286
$ws.Range("N51").Value = -10066.286
$ws.Range("H80").Value = 833.3333
$ws.Range("I80").Value = 833.3333
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2499.9999
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1501.9999
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 833.3333
$ws.Range("I83").Value = 833.3333
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 7499.9997
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -2507.9997
$ws.Range("N83").ClearContents()
$ws.Range("H88").Value = 12508819
$ws.Range("J88").Value = 10375.615
$ws.Range("L88").Value = 10375.615
$ws.Range("N88").Value = -11187.615
$ws.Range("H91").Value = 12508819
$ws.Range("J91").Value = 10375.615
$ws.Range("L91").Value = 10375.615
$ws.Range("N91").Value = -13183.615
$ws.Range("H101").Value = 2767.3635
$ws.Range("I101").Value = 363.75
$ws.Range("J101").Value = 4140.857
$ws.Range("K101").Value = 1091.25
$ws.Range("L101").Value = 12422.571
$ws.Range("M101").Value = 530.75
$ws.Range("N101").Value = -15666.571
$ws.Range("H107").Value = 266.875
$ws.Range("J107").Value = 251.28572
$ws.Range("L107").Value = 251.28572
$ws.Range("N107").Value = -4091.28572

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 674016.4399999999
$ws.Range("I2").Value = 833982.5600000001
$ws.Range("J2").Value = 2158.6
$ws.Range("K2").Value = 833982.5600000001
$ws.Range("L2").Value = 2158.6
$ws.Range("M2").Value = -833869.5600000001
$ws.Range("N2").Value = -2384.6
$ws.Range("H45").Value = 5666.3335
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 5666.3335
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 5666.3335
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -6420.3335
$ws.Range("H116").Value = 674016.4399999999
$ws.Range("I116").Value = 833982.5600000001
$ws.Range("J116").Value = 2158.6
$ws.Range("K116").Value = 833982.5600000001
$ws.Range("L116").Value = 2158.6
$ws.Range("M116").Value = -831688.5600000001
$ws.Range("N116").Value = -6746.6

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 674016.4399999999
$ws.Range("I3").Value = 833982.5600000001
$ws.Range("J3").Value = 2158.6
$ws.Range("K3").Value = 833982.5600000001
$ws.Range("L3").Value = 2158.6
$ws.Range("M3").Value = -833868.5600000001
$ws.Range("N3").Value = -2386.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1299
$ws.Range("I22").Value = 299.5
$ws.Range("K22").Value = 299.5
$ws.Range("M22").Value = 50.5
$ws.Range("H58").Value = 695049.4399999999
$ws.Range("I58").Value = 911721.25
$ws.Range("K58").Value = 911721.25
$ws.Range("M58").Value = -911518.25
$ws.Range("H105").Value = 3789245.2
$ws.Range("I105").Value = 11364236
$ws.Range("K105").Value = 11364236
$ws.Range("M105").Value = -11362489
$ws.Range("H122").Value = 3507.258
$ws.Range("I122").Value = 2085.05
$ws.Range("K122").Value = 6255.150000000001
$ws.Range("M122").Value = -3805.150000000001
$ws.Range("H132").Value = 7356.4473
$ws.Range("I132").Value = 8687.645500000001
$ws.Range("J132").Value = 1461.1428
$ws.Range("K132").Value = 26062.9365
$ws.Range("L132").Value = 4383.428400000001
$ws.Range("M132").Value = -23532.9365
$ws.Range("N132").Value = -9443.428400000001
$ws.Range("H134").Value = 1979.0526
$ws.Range("I134").Value = 1979.0526
$ws.Range("K134").Value = 5937.1578
$ws.Range("M134").Value = -3402.1578
$ws.Range("H136").Value = 695049.4399999999
$ws.Range("I136").Value = 911721.25
$ws.Range("K136").Value = 2735163.75
$ws.Range("M136").Value = -2732613.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 142.77777
$ws.Range("J12").Value = 105
$ws.Range("L12").Value = 315
$ws.Range("N12").Value = -661
$ws.Range("H26").Value = 111.44444
$ws.Range("I26").Value = 69.85714
$ws.Range("K26").Value = 209.57142
$ws.Range("M26").Value = 78.42858000000001
$ws.Range("H75").Value = 8147.1665
$ws.Range("I75").Value = 3074
$ws.Range("J75").Value = 9161.799999999999
$ws.Range("K75").Value = 9222
$ws.Range("L75").Value = 27485.4
$ws.Range("M75").Value = -8224
$ws.Range("N75").Value = -29481.4
$ws.Range("H78").Value = 8147.1665
$ws.Range("I78").Value = 3074
$ws.Range("J78").Value = 9161.799999999999
$ws.Range("K78").Value = 27666
$ws.Range("L78").Value = 82456.2
$ws.Range("M78").Value = -22674
$ws.Range("N78").Value = -92440.2
$ws.Range("H107").Value = 199990
$ws.Range("J107").Value = 199990
$ws.Range("L107").Value = 599970
$ws.Range("N107").Value = -603810
$ws.Range("H113").Value = 680.2941
$ws.Range("I113").Value = 719.2
$ws.Range("J113").Value = 664.0833
$ws.Range("K113").Value = 2157.6
$ws.Range("L113").Value = 1992.2499
$ws.Range("M113").Value = 12.39999999999964
$ws.Range("N113").Value = -6332.2499
$ws.Range("H122").Value = 1139.3
$ws.Range("J122").Value = 1287.7142
$ws.Range("L122").Value = 11589.4278
$ws.Range("N122").Value = -16489.4278
$ws.Range("H132").Value = 6489.5557
$ws.Range("I132").Value = 1085.2
$ws.Range("J132").Value = 13245
$ws.Range("K132").Value = 9766.800000000001
$ws.Range("L132").Value = 119205
$ws.Range("M132").Value = -7236.800000000001
$ws.Range("N132").Value = -124265

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 526.4761999999999
$ws.Range("I97").Value = 421.625
$ws.Range("J97").Value = 862
$ws.Range("K97").Value = 421.625
$ws.Range("L97").Value = 862
$ws.Range("M97").Value = 74.375
$ws.Range("N97").Value = -1854
$ws.Range("H107").Value = 9524169
$ws.Range("I107").Value = 11905112
$ws.Range("K107").Value = 11905112
$ws.Range("M107").Value = -11903192
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 526750.8
$ws.Range("I132").Value = 146448.22
$ws.Range("K132").Value = 439344.66
$ws.Range("M132").Value = -436814.66

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 66667770
$ws.Range("I40").Value = 1333
$ws.Range("K40").Value = 1333
$ws.Range("M40").Value = -1197
$ws.Range("H122").Value = 6874.5356
$ws.Range("I122").Value = 4139.2
$ws.Range("K122").Value = 12417.6
$ws.Range("M122").Value = -9967.599999999999
$ws.Range("H136").Value = 3880
$ws.Range("I136").Value = 2299.6667
$ws.Range("J136").Value = 6250.5
$ws.Range("K136").Value = 6899.000100000001
$ws.Range("L136").Value = 18751.5
$ws.Range("M136").Value = -4349.000100000001
$ws.Range("N136").Value = -23851.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 80721.125
$ws.Range("J80").Value = 80721.125
$ws.Range("L80").Value = 80721.125
$ws.Range("N80").Value = -82717.125
$ws.Range("H83").Value = 80721.125
$ws.Range("J83").Value = 80721.125
$ws.Range("L83").Value = 242163.375
$ws.Range("N83").Value = -252147.375
$ws.Range("H100").Value = 751706.25
$ws.Range("I100").Value = 959402.5600000001
$ws.Range("J100").Value = 3999.6
$ws.Range("K100").Value = 1918805.12
$ws.Range("L100").Value = 7999.2
$ws.Range("M100").Value = -1918264.12
$ws.Range("N100").Value = -9081.200000000001
$ws.Range("H132").Value = 16132291
$ws.Range("I132").Value = 1023.5833
$ws.Range("K132").Value = 3070.7499
$ws.Range("M132").Value = -540.7498999999998
$ws.Range("H136").Value = 10016.75
$ws.Range("I136").Value = 5098.643
$ws.Range("J136").Value = 11393.82
$ws.Range("K136").Value = 15295.929
$ws.Range("L136").Value = 34181.46
$ws.Range("M136").Value = -12745.929
$ws.Range("N136").Value = -39281.46

Write-Host "Applied all Cactuar_Profits updates"